$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVT")

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Grid Line"
$ws.Range("C9").Value = "Update grid lines"
$ws.Range("D9").Value = "1. Go to formatting pane`n2. Switch toggle of 'Grid Lines' to 'OFF'`n3. Switch toggle of 'Grid Lines' to 'ON'`n4. In Grid Lines field,  Switch toggle of 'X Axis' to 'OFF'`n5. In Grid Lines field,  Switch toggle of 'Y Axis' to 'OFF'"
$ws.Range("E9").Value = "Grid Lines will appear/disappear as per the toggle "

$ws.Range("D9").WrapText = $true
$ws.Range("E9").WrapText = $true
$ws.Range("C9").WrapText = $true

$ws.Rows.Item(9).RowHeight = 105

$ws.Range("E9").Select()
